# GSP469: rebuild the Overview sheet's bakery sales table
#  - header row gains a new "# Sold" column (D)
#  - header labels / day names get capitalized
#  - "pie" product data is replaced by distinct baked-goods per weekday
#  - the duplicate trailing "monday / pie" row (old row 9) is removed
#  - a "# Sold" quantity is added for each day (column D)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)   # "Overview"

# --- Header row ---------------------------------------------------------
$ws.Cells.Item(1,1).Value = "Day"
$ws.Cells.Item(1,2).Value = "Product"
$ws.Cells.Item(1,3).Value = "Price"
$ws.Cells.Item(1,4).Value = "# Sold"

# --- Data rows: Day, Product, Price, # Sold -----------------------------
$data = @(
    @("Monday",    "Croissant", 2.5, 120),
    @("Tuesday",   "Baguette",  1.5, 90),
    @("Wednesday", "Muffin",    2.0, 75),
    @("Thursday",  "Donut",     1.8, 130),
    @("Friday",    "Cupcake",   2.2, 110),
    @("Saturday",  "Bagel",     1.7, 95),
    @("Sunday",    "Brownie",   2.8, 100)
)

$row = 2
foreach ($rec in $data) {
    $ws.Cells.Item($row,1).Value = $rec[0]
    $ws.Cells.Item($row,2).Value = $rec[1]
    $ws.Cells.Item($row,3).Value = $rec[2]
    $ws.Cells.Item($row,4).Value = $rec[3]
    $row = $row + 1
}

# --- Remove the old duplicate 9th row (leftover "Monday / pie" entry) ---
$used = $ws.UsedRange
if ($used.Rows.Count -ge 9) {
    $ws.Rows(9).Delete()
}

# --- Match column D's cell style/format to the rest of the table (col C) -
$ws.Range("C1").Copy()
$ws.Range("D1:D8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false
